# Ajout log.addSTEPLOOP + corrections
# - A3 changes from "AD.SEC.002.FON.01" to "RO.ACT.001"
# - A4 gains "RO.ACT.003"
# - A5 gains "RO.ACT.004"
# - A6 (new row) gains "AD.SEC.002.FON.01"
# - D13 gains the new lookup value "RO.ACT.005"
# - Selection moves from A5 to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

$ws.Range("A3").Value = "RO.ACT.001"
$ws.Range("A4").Value = "RO.ACT.003"
$ws.Range("A5").Value = "RO.ACT.004"
$ws.Range("A6").Value = "AD.SEC.002.FON.01"
$ws.Range("D13").Value = "RO.ACT.005"

$ws.Range("B8").Select() | Out-Null
